$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("C36")
$rng.Borders.Item(10).LineStyle = 1
$rng.Borders.Item(10).Weight = 2
